$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 8226166.5
$ws.Range("I86").Value = 2605514
$ws.Range("K86").Value = 2605514
$ws.Range("M86").Value = -2604391
$ws.Range("H89").Value = 8226166.5
$ws.Range("I89").Value = 2605514
$ws.Range("K89").Value = 13027570
$ws.Range("M89").Value = -13021954
$ws.Range("H116").Value = 240390.33
$ws.Range("J116").Value = 298992.1
$ws.Range("L116").Value = 298992.1
$ws.Range("N116").Value = -305876.1
$ws.Range("H125","I125","J125","K125","L125","M125","N125","H126","I126","J126","K126","L126","H127","I127","J127","K127","L127","M127","H128","I128","J128","K128","L128","H129","I129","J129","K129","L129","M129","N129","H130","I130","J130","K130","L130","H131","I131","J131","K131","L131","M131","N131","H132","I132","J132","K132","L132","M132","N132","H133","I133","J133","K133","L133","N133","H134","I134","J134","K134","L134","H135","I135","J135","K135","L135","M135","N135","H136","I136","J136","K136","L136","N136","H137","I137","J137","K137","L137","M137","N137","H138","I138","J138","K138","L138","M138","N138","H139","I139","J139","K139","L139","M139","N139","H140","I140","J140","K140","L140","N140","H141","I141","J141","K141","L141","M141","N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 20000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 20000
$ws.Range("K39").Value = 0
$ws.Range("M39").Value = 20000
$ws.Range("N39").Value = -21040
$ws.Range("H122").Value = 3889.0588
$ws.Range("I122").Value = 3340.125
$ws.Range("J122").Value = 4377
$ws.Range("K122").Value = 10020.375
$ws.Range("L122").Value = 13131
$ws.Range("M122").Value = -7570.375
$ws.Range("N122").Value = -18031
$ws.Range("H132").Value = 1821676.2
$ws.Range("I132").Value = 1821676.2
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5465028.6
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -5462498.6
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("N133").Value = 0
$ws.Range("H137").Value = 85000
$ws.Range("J137").Value = 85000
$ws.Range("L137").Value = 85000
$ws.Range("N137").Value = -95200
$ws.Range("L39","M132","L133").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 26260
$ws.Range("I96").Value = 26260
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 26260
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = -23514
$ws.Range("H117").Value = 100742
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 100742
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 100742
$ws.Range("N117").Value = -109920
$ws.Range("H118").Value = 50000
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 50000
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 50000
$ws.Range("N118").Value = -53314
$ws.Range("H119").Value = 100761
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 100761
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 100761
$ws.Range("N119").Value = -110437
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("H122").Value = 129000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 129000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 129000
$ws.Range("N122").Value = -138800
$ws.Range("H123").Value = 75815
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 75815
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 75815
$ws.Range("N123").Value = -85615
$ws.Range("H124").Value = 55000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 55000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 55000
$ws.Range("N124").Value = -64820
$ws.Range("H125").Value = 78000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 78000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 78000
$ws.Range("N125").Value = -87840
$ws.Range("H126").Value = 44000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 44000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 44000
$ws.Range("N126").Value = -53880
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H128").Value = 16333
$ws.Range("I128").Value = 16333
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 48999
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -46509
$ws.Range("H129").Value = 79000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 79000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 79000
$ws.Range("N129").Value = -89000
$ws.Range("H130").Value = 233593.33
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 233593.33
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 233593.33
$ws.Range("N130").Value = -243633.33
$ws.Range("H131").Value = 64666.168
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 64666.168
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 64666.168
$ws.Range("N131").Value = -74746.16800000001
$ws.Range("H132").Value = 499999
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 499999
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 499999
$ws.Range("N132").Value = -510119
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H134").Value = 3445728
$ws.Range("I134").Value = 5653326.5
$ws.Range("J134").Value = 256974.72
$ws.Range("K134").Value = 16959979.5
$ws.Range("L134").Value = 770924.16
$ws.Range("M134").Value = -16957444.5
$ws.Range("N134").Value = -775994.16
$ws.Range("H135").Value = 60000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 60000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 99000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 99000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 99000
$ws.Range("N139").Value = -109280
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M96").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 232.23077
$ws.Range("I7").Value = 225.5
$ws.Range("J7").Value = 235.22223
$ws.Range("K7").Value = 225.5
$ws.Range("L7").Value = 235.22223
$ws.Range("M7").Value = -112.5
$ws.Range("N7").Value = -461.22223
$ws.Range("H58").Value = 364900.2
$ws.Range("I58").Value = 537751.9399999999
$ws.Range("J58").Value = 3482.818
$ws.Range("K58").Value = 537751.9399999999
$ws.Range("L58").Value = 3482.818
$ws.Range("M58").Value = -537548.9399999999
$ws.Range("N58").Value = -3888.818
$ws.Range("H59").Value = 119996
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 119996
$ws.Range("K59").Value = 0
$ws.Range("M59").Value = 119996
$ws.Range("N59").Value = -122286
$ws.Range("H62").Value = 6803.625
$ws.Range("I62").Value = 2512.5
$ws.Range("J62").Value = 8234
$ws.Range("K62").Value = 2512.5
$ws.Range("L62").Value = 8234
$ws.Range("M62").Value = -1888.5
$ws.Range("N62").Value = -9482
$ws.Range("H65").Value = 6803.625
$ws.Range("I65").Value = 2512.5
$ws.Range("J65").Value = 8234
$ws.Range("K65").Value = 12562.5
$ws.Range("L65").Value = 41170
$ws.Range("M65").Value = -9442.5
$ws.Range("N65").Value = -47410
$ws.Range("H68").Value = 79088.39999999999
$ws.Range("J68").Value = 96360.5
$ws.Range("L68").Value = 96360.5
$ws.Range("N68").Value = -97858.5
$ws.Range("H71").Value = 79088.39999999999
$ws.Range("J71").Value = 96360.5
$ws.Range("L71").Value = 289081.5
$ws.Range("N71").Value = -296569.5
$ws.Range("H94").Value = 2336.8667
$ws.Range("I94").Value = 1739.75
$ws.Range("K94").Value = 1739.75
$ws.Range("M94").Value = -1288.75
$ws.Range("H99").Value = 3684.875
$ws.Range("I99").Value = 2750
$ws.Range("K99").Value = 2750
$ws.Range("M99").Value = -1252
$ws.Range("H126").Value = 3684.875
$ws.Range("I126").Value = 2750
$ws.Range("K126").Value = 8250
$ws.Range("M126").Value = -5780
$ws.Range("H132").Value = 2489287
$ws.Range("J132").Value = 64061.5
$ws.Range("L132").Value = 192184.5
$ws.Range("N132").Value = -197244.5
$ws.Range("H134").Value = 13284.091
$ws.Range("I134").Value = 15569.889
$ws.Range("J134").Value = 2998
$ws.Range("K134").Value = 46709.667
$ws.Range("L134").Value = 8994
$ws.Range("M134").Value = -44174.667
$ws.Range("N134").Value = -14064
$ws.Range("H136").Value = 364900.2
$ws.Range("I136").Value = 537751.9399999999
$ws.Range("J136").Value = 3482.818
$ws.Range("K136").Value = 1613255.82
$ws.Range("L136").Value = 10448.454
$ws.Range("M136").Value = -1610705.82
$ws.Range("N136").Value = -15548.454
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("N140").Value = 0
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("N141").Value = 0
$ws.Range("L59","M137","L140","L141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5951
$ws.Range("I68").Value = 1732.2858
$ws.Range("K68").Value = 5196.857400000001
$ws.Range("M68").Value = -4385.857400000001
$ws.Range("H71").Value = 5951
$ws.Range("I71").Value = 1732.2858
$ws.Range("K71").Value = 15590.5722
$ws.Range("M71").Value = -11534.5722
$ws.Range("H107").Value = 908.4286
$ws.Range("J107").Value = 873.8333
$ws.Range("L107").Value = 2621.4999
$ws.Range("N107").Value = -6461.4999
$ws.Range("H120","I120","J120","K120","L120","M120","N120","H121","I121","J121","K121","L121","M121","N121","H122","I122","J122","K122","L122","M122","N122","H123","I123","J123","K123","L123","M123","N123","H124","I124","J124","K124","L124","M124","N124","H125","I125","J125","K125","L125","M125","H126","I126","J126","K126","L126","M126","N126","H127","I127","J127","K127","L127","N127","H128","I128","J128","K128","L128","M128","H129","I129","J129","K129","L129","M129","N129","H130","I130","J130","K130","L130","M130","N130","H131","I131","J131","K131","L131","M131","N131","H132","I132","J132","K132","L132","M132","N132","H133","I133","J133","K133","L133","M133","H134","I134","J134","K134","L134","M134","H136","I136","J136","K136","L136","M136","N136","H137","I137","J137","K137","L137","M137","N137","H138","I138","J138","K138","L138","M138","N138","H139","I139","J139","K139","L139","M139","N139","H140","I140","J140","K140","L140","M140","N140","H141","I141","J141","K141","L141","M141","N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 22182.584
$ws.Range("J57").Value = 28810.111
$ws.Range("L57").Value = 28810.111
$ws.Range("N57").Value = -30450.111
$ws.Range("H132").Value = 526713
$ws.Range("I132").Value = 806705.1
$ws.Range("K132").Value = 2420115.3
$ws.Range("M132").Value = -2417585.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3955.524
$ws.Range("I40").Value = 3653.35
$ws.Range("K40").Value = 3653.35
$ws.Range("M40").Value = -3517.35
$ws.Range("H46").Value = 3246.2307
$ws.Range("J46").Value = 3745.2
$ws.Range("L46").Value = 3745.2
$ws.Range("N46").Value = -4121.2
$ws.Range("H63").Value = 76272
$ws.Range("J63").Value = 75003.664
$ws.Range("L63").Value = 75003.664
$ws.Range("N63").Value = -76501.664
$ws.Range("H66").Value = 76272
$ws.Range("J66").Value = 75003.664
$ws.Range("L66").Value = 225010.992
$ws.Range("N66").Value = -232498.992
$ws.Range("H122").Value = 5377.3887
$ws.Range("I122").Value = 5558.8335
$ws.Range("K122").Value = 16676.5005
$ws.Range("M122").Value = -14226.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5442225
$ws.Range("I132").Value = 6494158.5
$ws.Range("K132").Value = 19482475.5
$ws.Range("M132").Value = -19479945.5
